## Generate Report for Handback
## - Update status text from "Ready for handoff" to "Handed back: in sync with en-US"
##   on the Overview sheet (zh-cn/de-de status columns) and on each language sheet's
##   "Status" column.
## - Populate "Latest Target File" (col I) with the source file name + a hyperlink to
##   the same target as column A, for each localized-file row.
## - Populate "Latest Handback File" (col J) with the same text as "Latest Handoff
##   File" (col G) for each row, for each language sheet.
## - Stamp "Latest Handback DateTime" (col K): zh-cn rows -> 2016-11-29 02:54:13,
##   de-de rows -> 2016-11-29 02:54:31.
## - Widen a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for rows 2-6
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
for ($r = 2; $r -le 6; $r++) {
    $overview.Cells.Item($r, 5).Value = $statusText
    $overview.Cells.Item($r, 6).Value = $statusText
}
$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# Per-language sheets: zh-cn, de-de
# ---------------------------------------------------------------------------
$languages = @(
    @{ Name = "zh-cn"; HandbackDateTime = "2016-11-29 02:54:13" },
    @{ Name = "de-de"; HandbackDateTime = "2016-11-29 02:54:31" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Grab the existing column-A hyperlink targets (row -> URL) before we start
    # adding new ones, keyed by row number.
    $targets = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -match '\$A\$(\d+)') {
            $targets[[int]$matches[1]] = $hl.Address
        }
    }

    for ($r = 2; $r -le 6; $r++) {
        # Status column (C)
        $ws.Cells.Item($r, 3).Value = $statusText

        # Latest Target File (I) = Source File Name (A), with a hyperlink matching
        # the one already on column A.
        $sourceName = $ws.Cells.Item($r, 1).Value2
        $ws.Cells.Item($r, 9).Value = $sourceName
        $target = $targets[$r]
        if ($target) {
            $ws.Hyperlinks.Add($ws.Cells.Item($r, 9), $target, "", "", $sourceName) | Out-Null
        }

        # Latest Handback File (J) = Latest Handoff File (G)
        $handoffName = $ws.Cells.Item($r, 7).Value2
        $ws.Cells.Item($r, 10).Value = $handoffName

        # Latest Handback DateTime (K)
        $ws.Cells.Item($r, 11).Value = $lang.HandbackDateTime
    }

    $ws.Columns.Item(3).ColumnWidth = 29.14
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}
